$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Add new feature (column H) / notes (column I) annotations for rows 152-271 ---
$ws.Cells.Item(152, 8).Value = 'FT0227'
$ws.Cells.Item(152, 9).Value = 'Peak shifts a lot'
$ws.Cells.Item(153, 8).Value = 'FT0253'
$ws.Cells.Item(154, 8).Value = 'NA'
$ws.Cells.Item(155, 8).Value = 'NA'
$ws.Cells.Item(156, 9).Value = 'Unpicked'
$ws.Cells.Item(157, 8).Value = 'NA'
$ws.Cells.Item(158, 8).Value = 'NA'
$ws.Cells.Item(159, 8).Value = 'NA'
$ws.Cells.Item(160, 8).Value = 'NA'
$ws.Cells.Item(161, 8).Value = 'NA'
$ws.Cells.Item(162, 9).Value = 'Unpicked'
$ws.Cells.Item(163, 8).Value = 'NA'
$ws.Cells.Item(164, 8).Value = 'FT1201'
$ws.Cells.Item(165, 8).Value = 'NA'
$ws.Cells.Item(166, 8).Value = 'NA'
$ws.Cells.Item(167, 8).Value = 'FT1376'
$ws.Cells.Item(168, 8).Value = 'FT1391'
$ws.Cells.Item(169, 8).Value = 'FT1397'
$ws.Cells.Item(170, 9).Value = 'Unpicked'
$ws.Cells.Item(171, 8).Value = 'NA'
$ws.Cells.Item(172, 8).Value = 'FT1452'
$ws.Cells.Item(173, 8).Value = 'NA'
$ws.Cells.Item(174, 8).Value = 'NA'
$ws.Cells.Item(175, 8).Value = 'NA'
$ws.Cells.Item(176, 8).Value = 'FT1557'
$ws.Cells.Item(177, 8).Value = 'FT1616'
$ws.Cells.Item(178, 8).Value = 'FT1664'
$ws.Cells.Item(179, 8).Value = 'NA'
$ws.Cells.Item(180, 8).Value = 'FT1666'
$ws.Cells.Item(181, 8).Value = 'FT1732'
$ws.Cells.Item(182, 8).Value = 'NA'
$ws.Cells.Item(183, 8).Value = 'NA'
$ws.Cells.Item(184, 8).Value = 'NA'
$ws.Cells.Item(185, 8).Value = 'NA'
$ws.Cells.Item(186, 8).Value = 'FT1861'
$ws.Cells.Item(187, 8).Value = 'NA'
$ws.Cells.Item(188, 8).Value = 'FT1918'
$ws.Cells.Item(189, 8).Value = 'NA'
$ws.Cells.Item(190, 8).Value = 'NA'
$ws.Cells.Item(191, 8).Value = 'NA'
$ws.Cells.Item(192, 8).Value = 'NA'
$ws.Cells.Item(193, 8).Value = 'FT2039'
$ws.Cells.Item(194, 8).Value = 'NA'
$ws.Cells.Item(195, 8).Value = 'NA'
$ws.Cells.Item(196, 8).Value = 'NA'
$ws.Cells.Item(196, 9).Value = 'Maybe the gross peak at 9.5, super low quality'
$ws.Cells.Item(197, 8).Value = 'NA'
$ws.Cells.Item(198, 8).Value = 'NA'
$ws.Cells.Item(198, 9).Value = 'Peak looks terrible'
$ws.Cells.Item(199, 8).Value = 'NA'
$ws.Cells.Item(200, 8).Value = 'NA'
$ws.Cells.Item(201, 8).Value = 'NA'
$ws.Cells.Item(202, 8).Value = 'NA'
$ws.Cells.Item(203, 8).Value = 'NA'
$ws.Cells.Item(204, 8).Value = 'NA'
$ws.Cells.Item(205, 8).Value = 'NA'
$ws.Cells.Item(206, 8).Value = 'NA'
$ws.Cells.Item(207, 8).Value = 'NA'
$ws.Cells.Item(208, 8).Value = 'NA'
$ws.Cells.Item(209, 8).Value = 'NA'
$ws.Cells.Item(210, 8).Value = 'FT2435'
$ws.Cells.Item(210, 9).Value = 'Peak looks terrible'
$ws.Cells.Item(211, 8).Value = 'NA'
$ws.Cells.Item(212, 8).Value = 'NA'
$ws.Cells.Item(213, 8).Value = 'NA'
$ws.Cells.Item(214, 8).Value = 'NA'
$ws.Cells.Item(215, 8).Value = 'FT2626'
$ws.Cells.Item(215, 9).Value = 'Maybe, standards look terrible'
$ws.Cells.Item(216, 8).Value = 'NA'
$ws.Cells.Item(217, 8).Value = 'NA'
$ws.Cells.Item(218, 8).Value = 'NA'
$ws.Cells.Item(219, 8).Value = 'NA'
$ws.Cells.Item(220, 8).Value = 'NA'
$ws.Cells.Item(221, 8).Value = 'NA'
$ws.Cells.Item(222, 8).Value = 'NA'
$ws.Cells.Item(223, 8).Value = 'NA'
$ws.Cells.Item(224, 8).Value = 'NA'
$ws.Cells.Item(225, 8).Value = 'NA'
$ws.Cells.Item(226, 8).Value = 'NA'
$ws.Cells.Item(227, 8).Value = 'NA'
$ws.Cells.Item(228, 8).Value = 'FT3204'
$ws.Cells.Item(228, 9).Value = 'Meh peak'
$ws.Cells.Item(229, 8).Value = 'NA'
$ws.Cells.Item(229, 9).Value = 'Looks fine in stans but not present in samps'
$ws.Cells.Item(230, 8).Value = 'NA'
$ws.Cells.Item(231, 8).Value = 'NA'
$ws.Cells.Item(232, 8).Value = 'NA'
$ws.Cells.Item(233, 8).Value = 'NA'
$ws.Cells.Item(234, 8).Value = 'NA'
$ws.Cells.Item(235, 8).Value = 'NA'
$ws.Cells.Item(236, 8).Value = 'NA'
$ws.Cells.Item(237, 8).Value = 'NA'
$ws.Cells.Item(238, 8).Value = 'NA'
$ws.Cells.Item(239, 8).Value = 'NA'
$ws.Cells.Item(240, 8).Value = 'NA'
$ws.Cells.Item(241, 8).Value = 'NA'
$ws.Cells.Item(242, 8).Value = 'FT3571'
$ws.Cells.Item(243, 8).Value = 'NA'
$ws.Cells.Item(244, 8).Value = 'NA'
$ws.Cells.Item(245, 8).Value = 'NA'
$ws.Cells.Item(246, 8).Value = 'NA'
$ws.Cells.Item(247, 8).Value = 'NA'
$ws.Cells.Item(248, 8).Value = 'NA'
$ws.Cells.Item(249, 8).Value = 'NA'
$ws.Cells.Item(250, 8).Value = 'NA'
$ws.Cells.Item(251, 8).Value = 'NA'
$ws.Cells.Item(252, 8).Value = 'NA'
$ws.Cells.Item(253, 8).Value = 'NA'
$ws.Cells.Item(254, 8).Value = 'NA'
$ws.Cells.Item(255, 8).Value = 'NA'
$ws.Cells.Item(256, 8).Value = 'NA'
$ws.Cells.Item(257, 8).Value = 'NA'
$ws.Cells.Item(258, 8).Value = 'NA'
$ws.Cells.Item(259, 8).Value = 'NA'
$ws.Cells.Item(260, 8).Value = 'NA'
$ws.Cells.Item(261, 8).Value = 'NA'
$ws.Cells.Item(262, 8).Value = 'NA'
$ws.Cells.Item(263, 8).Value = 'NA'
$ws.Cells.Item(264, 8).Value = 'NA'
$ws.Cells.Item(265, 8).Value = 'NA'
$ws.Cells.Item(266, 8).Value = 'NA'
$ws.Cells.Item(267, 8).Value = 'NA'
$ws.Cells.Item(268, 8).Value = 'NA'
$ws.Cells.Item(269, 8).Value = 'NA'
$ws.Cells.Item(270, 8).Value = 'NA'
$ws.Cells.Item(271, 8).Value = 'NA'

# --- Update the view state to reflect where the user was working ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 196
$win.ScrollColumn = 1
$ws.Range("E200").Select()
